# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @(
    @{Row=2; B=0.1554434735375247; C=0.05231270169004087; D=0.7127328510149897; E=0.4998867070740569; G=1.420375733316612},
    @{Row=3; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=4; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054},
    @{Row=5; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729},
    @{Row=6; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054},
    @{Row=7; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126},
    @{Row=8; B=0.02258322285507441; C=0.004309184025731883; D=0.7127328510149897; E=0.4998867070740569; G=1.239511964969853},
    @{Row=9; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=10; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=11; B=1.505614041169197; C=0.3375848360084654; D=0.1529057820181812; E=0.4998867070740569; G=2.495991366269901},
    @{Row=12; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=13; B=1.505614041169197; C=0.3375848360084654; D=0.1529057820181812; E=0.4998867070740569; G=2.495991366269901},
    @{Row=14; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126},
    @{Row=15; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729},
    @{Row=16; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=17; B=1.505614041169197; C=1.65323645889881; D=0.1529057820181812; E=6.48142807727062; G=9.793184359356808},
    @{Row=18; B=0.7287194209349384; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=3.594575437922795},
    @{Row=19; B=1.505614041169197; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=6.741336633845642},
    @{Row=20; B=0.06328177979961902; C=0.3375848360084654; D=16.98373111632243; E=6.48142807727062; G=23.86602580940113},
    @{Row=21; B=0.1554434735375247; C=0.3375848360084654; D=16.98373111632243; E=6.48142807727062; G=23.95818750313904},
    @{Row=22; B=1.505614041169197; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=6.741336633845642},
    @{Row=23; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054},
    @{Row=24; B=0.06328177979961902; C=1.65323645889881; D=0.1529057820181812; E=6.48142807727062; G=8.35085209798723},
    @{Row=25; B=3.182878228561681; C=1.65323645889881; D=0.1529057820181812; E=0.4998867070740569; G=5.488907176552729},
    @{Row=26; B=0.7287194209349384; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=3.594575437922795},
    @{Row=27; B=3.182878228561681; C=1.65323645889881; D=157.8057217802531; E=6.48142807727062; G=169.1232645449842},
    @{Row=28; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126},
    @{Row=29; B=0.3464964993005633; C=9.226618575922256; D=157.8057217802531; E=6.48142807727062; G=173.8602649327466},
    @{Row=30; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126},
    @{Row=31; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126},
    @{Row=32; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054},
    @{Row=33; B=3.182878228561681; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=6.048734245549538},
    @{Row=34; B=1.505614041169197; C=1.65323645889881; D=0.7127328510149897; E=0.4998867070740569; G=4.371470058157054},
    @{Row=35; B=3.182878228561681; C=1.65323645889881; D=3.082599426703578; E=0.4998867070740569; G=8.418600821238126}
)

foreach ($u in $updates) {
    $row = $u.Row
    $ws.Range("B$row").Value = $u.B
    $ws.Range("C$row").Value = $u.C
    $ws.Range("D$row").Value = $u.D
    $ws.Range("E$row").Value = $u.E
    $ws.Range("G$row").Value = $u.G
}
